$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has an AutoFilter over A1:E220 that was previously showing only
# the "application_{train|test}.csv" rows (column B). Re-apply the filter
# to show "bureau.csv" rows instead - the author has moved on to looking at
# the bureau files. Passing Criteria1 as an array (rather than a bare
# string) makes Excel record this as a discrete value-list pick
# (<filters><filter val="bureau.csv"/></filters>) instead of a comparison
# <customFilters>, matching what a normal dropdown checkbox selection
# produces. Field 2 = column B (second column of the filtered range).
$rng = $ws.Range("A1:E220")
[void]$rng.AutoFilter(2, @("bureau.csv"))

# Row 128 (LIVE_ID 129) has a leftover "bureau.csv" value sitting in column
# B that doesn't belong - every other bureau.csv row leaves column B blank.
# Clear it now that the filter has already used its original value to
# decide the row belongs in the visible "bureau.csv" set.
[void]$ws.Range("B128").ClearContents()

# Move the cursor/selection to where the user was now working.
[void]$ws.Range("D131").Select()

Write-Host "done"
